$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New data rows to append (Fecha serial, Pruebas Realizadas, Pruebas Positivas,
# Clinicamente Estables, Clinicamente Graves, Cuidados Intensivos)
$newData = @(
    @(44139, 1205, 436, 404, 107, 19),
    @(44140, 906, 223, 392, 96, 19),
    @(44141, 969, 229, 354, 105, 21),
    @(44142, 1653, 465, 295, 92, 15)
)

foreach ($rowVals in $newData) {
    # Grow the table by one row (keeps table ref / autoFilter / dimension in sync)
    $newRow = $lo.ListRows.Add()
    $r = $lo.Range.Row + $lo.Range.Rows.Count - 1

    # Copy formatting from the row directly above so the new row matches the
    # existing date / centered-number styling used throughout the table.
    $ws.Range($ws.Cells.Item($r - 1, 1), $ws.Cells.Item($r - 1, 6)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 6)).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
}

$excel.CutCopyMode = 0

[void]$ws.Range("D243").Select()
